$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "58.797.49"
Set-TextValue "E2" "  -2.89%  "
Set-TextValue "D3" "2.726.00"
Set-TextValue "E3" "  -6.32%  "
Set-TextValue "E4" "  +0.02%  "
Set-TextValue "D5" "505.12"
Set-TextValue "E5" "  -4.18%  "
Set-TextValue "D6" "141.52"
Set-TextValue "E6" "  -0.22%  "
Set-TextValue "D7" "0.998"
Set-TextValue "E7" "  -0.14%  "
Set-TextValue "D8" "0.531"
Set-TextValue "E8" "  -3.52%  "
Set-TextValue "D9" "2.736.05"
Set-TextValue "E9" "  -6.17%  "
Set-TextValue "D10" "6.06"
Set-TextValue "E10" "  +2.62%  "
Set-TextValue "E11" "  -2.24%  "
Set-TextValue "D12" "0.349"
Set-TextValue "E12" "  -2.12%  "
Set-TextValue "E13" "  +1.13%  "
Set-TextValue "D14" "3.208.43"
Set-TextValue "E14" "  -6.23%  "
Set-TextValue "D15" "58.909.33"
Set-TextValue "E15" "  -2.79%  "
Set-TextValue "D16" "21.67"
Set-TextValue "E16" "  -4.43%  "
Set-TextValue "B17" "WrappedEther"
Set-TextValue "C17" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D17" "2.741.45"
Set-TextValue "E17" "  -5.70%  "
Set-TextValue "B18" "ShibaInu"
Set-TextValue "C18" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D18" "0.0000136"
Set-TextValue "E18" "  -3.28%  "
Set-TextValue "D19" "4.76"
Set-TextValue "E19" "  -3.71%  "
Set-TextValue "D20" "11.02"
Set-TextValue "E20" "  -4.56%  "
Set-TextValue "D21" "344.52"
Set-TextValue "E21" "  -4.66%  "
Set-TextValue "D22" "6.25"
Set-TextValue "E22" "  -4.93%  "
Set-TextValue "D23" "0.998"
Set-TextValue "E23" "  -0.23%  "
Set-TextValue "D24" "5.61"
Set-TextValue "E24" "  -0.48%  "
Set-TextValue "D25" "63.22"
Set-TextValue "E25" "  -0.30%  "
Set-TextValue "D26" "0.427"
Set-TextValue "E26" "  -4.79%  "
Set-TextValue "D27" "0.172"
Set-TextValue "E27" "  -4.51%  "
Set-TextValue "D28" "0.994"
Set-TextValue "E28" "  -0.34%  "
Set-TextValue "B29" "InternetComputer(DFINITY)"
Set-TextValue "C29" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D29" "7.51"
Set-TextValue "E29" "  -3.77%  "
Set-TextValue "B30" "PEPE"
Set-TextValue "C30" "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue "D30" "0.0₃0833"
Set-TextValue "E30" "  -2.07%  "
Set-TextValue "E31" "  -0.09%  "
Set-TextValue "B32" "PancakeSwap"
Set-TextValue "C32" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D32" "1.61"
Set-TextValue "E32" "  -3.49%  "
Set-TextValue "B33" "EthereumClassic"
Set-TextValue "C33" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D33" "19.19"
Set-TextValue "E33" "  -1.51%  "
Set-TextValue "D34" "150.16"
Set-TextValue "E34" "  -0.53%  "
Set-TextValue "B35" "Aptos"
Set-TextValue "C35" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D35" "5.42"
Set-TextValue "E35" "  -2.06%  "
Set-TextValue "B36" "NEARProtocol"
Set-TextValue "C36" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D36" "4.19"
Set-TextValue "E36" "  -2.98%  "
Set-TextValue "D37" "0.948"
Set-TextValue "E37" "  -4.05%  "
Set-TextValue "E38" "  -5.58%  "
Set-TextValue "D39" "36.02"
Set-TextValue "E39" "  -5.25%  "
Set-TextValue "D40" "1.39"
Set-TextValue "E40" "  -5.67%  "
Set-TextValue "D41" "3.54"
Set-TextValue "E41" "  -3.10%  "
Set-TextValue "D42" "2.184.89"
Set-TextValue "E42" "  -6.43%  "
Set-TextValue "B43" "FirstDigitalUSD"
Set-TextValue "C43" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D43" "0.996"
Set-TextValue "E43" "  -0.19%  "
Set-TextValue "B44" "Hedera"
Set-TextValue "C44" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D44" "0.0557"
Set-TextValue "E44" "  -2.56%  "
Set-TextValue "D45" "0.604"
Set-TextValue "E45" "  -6.50%  "
Set-TextValue "D46" "19.04"
Set-TextValue "E46" "  -8.30%  "
Set-TextValue "E47" "  +0.24%  "
Set-TextValue "D48" "4.76"
Set-TextValue "E48" "  -4.01%  "
Set-TextValue "D49" "0.0226"
Set-TextValue "E49" "  -3.25%  "
Set-TextValue "D50" "0.0887"
Set-TextValue "E50" "  -4.39%  "
Set-TextValue "D51" "18.05"
Set-TextValue "E51" "  -1.05%  "
